$d = $word.ActiveDocument

# 1) Title: "Call Logger v0.1.3 Instruction Manual" -> "...v0.1.4..."
$null = $d.Content.Find.Execute(
    "Call Logger v0.1.3 Instruction Manual", $true, $false, $false, $false, $false,
    $true, 1, $false, "Call Logger v0.1.4 Instruction Manual", 2)

# 2) Heading: merge "Advanced: " + "Outcome Column" into a single run's text
$null = $d.Content.Find.Execute(
    "Advanced: Outcome Column", $true, $false, $false, $false, $false,
    $true, 1, $false, "Advanced: Outcome Column", 2)
